$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 changes
$ws.Range("AI2").Value = 800

# Row 3 changes
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5
$ws.Range("AP3").Value = 1.88
$ws.Range("AQ3").Value = 1.98

# Row 4 changes
$ws.Range("Q4").Value = 1.93
$ws.Range("R4").Value = 1.93
